$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the title heading ---
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaBoldText = "Meta description"
$metaRestText = ": Discover the charming design and exciting features of Beellionaires Dream Drop. Play now for free and potentially win big with random multipliers and progressive jackpots."

$metaStart = $metaPara.Range.Start
$fullMetaRange = $d.Range($metaStart, $metaStart)
$fullMetaRange.Text = $metaBoldText + $metaRestText

$boldRange = $d.Range($metaStart, $metaStart + $metaBoldText.Length)
$boldRange.Bold = 1

# --- 2. Remove the trailing duplicate "Play Beellionaires Dream Drop for Free - Review" paragraph ---
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# --- 3. Replace the italic "meta description" paragraph at the end with the DALLE image prompt ---
$count2 = $d.Paragraphs.Count
$imgPara = $d.Paragraphs.Item($count2)
$imgStart = $imgPara.Range.Start
$imgEnd = $imgPara.Range.End

$newImgText = 'Create a feature image fitting "Beellionaires Dream Drop": - The image depicts a happy Maya warrior wearing glasses - The image is in a cartoon style For the feature image of "Beellionaires Dream Drop", DALLE can create an illustration of a Maya warrior celebrating amidst a beehive. The warrior is wearing glasses that give a nerdy yet fun look to the character. The bee theme of the game can be represented in the image by including bees or honeycombs in the background. The artwork should be colorful and vibrant to attract the attention of the players. The cartoon style of the image gives the game a playful and lighthearted feel, which matches the theme of the game. The format of the image can be square or rectangular to fit the aesthetics of the online casino or social media platform where it will be posted.'

$imgRange = $d.Range($imgStart, $imgEnd - 1)
$imgRange.Text = $newImgText
$imgRange.Italic = 1

Write-Output "done"
